$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update player names (column B)
$ws.Range("B2").Value = "Anthony Davis"
$ws.Range("B3").Value = "Nikola Jokic"
$ws.Range("B4").Value = "Josh Hart"
$ws.Range("B5").Value = "Aaron Gordon"
$ws.Range("B6").Value = "Bobby Portis"

# Update team names (column C)
$ws.Range("C2").Value = "LAL"
$ws.Range("C3").Value = "DEN"
$ws.Range("C4").Value = "NY"
$ws.Range("C5").Value = "DEN"
$ws.Range("C6").Value = "MIL"

# Update rank (column A)
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 3

# Update quantity (column D)
$ws.Range("D2").Value = 3
$ws.Range("D3").Value = 3
$ws.Range("D4").Value = 2
$ws.Range("D5").Value = 2
$ws.Range("D6").Value = 2
